# docs(wbs): refresh milestone and row execution statuses
#
# Bump the "Status Updated On" date (column I) for every data row of the
# WBS sheet from 2026-02-24 (Excel serial 46077) to 2026-02-25 (Excel
# serial 46078).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WBS")

for ($row = 2; $row -le 137; $row++) {
    $cell = $ws.Cells.Item($row, 9)  # Column I
    if ($cell.Value2 -eq 46077) {
        $cell.Value2 = 46078
    }
}
